# Add 2022-Q3 sheet data, matching commit "feat: add 2022-Q3 data"
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: force a cell to be stored as TEXT (not auto-converted to a
# number) and strip any inherited/implicit style so it ends up with the
# default (no style index) formatting.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ===========================================================================
# 1) Insert the new "2022-Q3" worksheet right after "总计" (i.e. before the
#    existing "2022-Q2" sheet), so tab order becomes:
#    总计, 2022-Q3, 2022-Q2, 2022-Q1
# ===========================================================================
$q2Sheet = $wb.Worksheets.Item(2)
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Reference sheet to copy header / column-A cell formatting from.
$refSheet = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Header row (row 1): B1:H1 -- copy style from the reference sheet header
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q3Sheet.Cells.Item(1, $col).Value = $headers[$col - 2]
    $refSheet.Cells.Item(1, $col).Copy()
    $q3Sheet.Cells.Item(1, $col).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# Data rows 2-7
# ---------------------------------------------------------------------
$rows = @(
    @("013357", "大摩沪港深精选混合C",       "1.53", "92.27", "5.83", "0.0892", 7),
    @("014114", "广发沪港深医药混合A",       "2.23", "93.52", "2.39", "0.0533", 10),
    @("013356", "大摩沪港深精选混合A",       "0.68", "92.27", "5.83", "0.0396", 7),
    @("014115", "广发沪港深医药混合C",       "1.04", "93.52", "2.39", "0.0249", 10),
    @("014214", "光大保德信核心资产混合A", "0.32", "88.01", "5.08", "0.0163", 4),
    @("014215", "光大保德信核心资产混合C", "0.03", "88.01", "5.08", "0.0015", 4)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    # Column A: 0-based index, numeric, styled like the reference sheet.
    $refSheet.Cells.Item(2, 1).Copy()
    $q3Sheet.Cells.Item($r, 1).PasteSpecial(-4122)
    $q3Sheet.Cells.Item($r, 1).Value = $i

    # Column B: fund code -- text (would otherwise parse as a number).
    Set-TextValue $q3Sheet.Cells.Item($r, 2) $data[0]

    # Column C: fund name -- plain text.
    $q3Sheet.Cells.Item($r, 3).Value = $data[1]

    # Columns D-G: decimal-looking values stored as TEXT in the source file.
    Set-TextValue $q3Sheet.Cells.Item($r, 4) $data[2]
    Set-TextValue $q3Sheet.Cells.Item($r, 5) $data[3]
    Set-TextValue $q3Sheet.Cells.Item($r, 6) $data[4]
    Set-TextValue $q3Sheet.Cells.Item($r, 7) $data[5]

    # Column H: rank -- numeric.
    $q3Sheet.Cells.Item($r, 8).Value = $data[6]
}

# ===========================================================================
# 2) Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q3,
#    pushing the existing 2022-Q2 / 2022-Q1 rows down to rows 3 / 4.
# ===========================================================================
$totalSheet = $wb.Worksheets.Item(1)
$totalSheet.Rows.Item(2).Insert()

# Clear the style the row-insert implicitly inherited on B2:D2.
$totalSheet.Range("B2:D2").Style = "Normal"

# Column A: numeric index, same style as the row below it (A3).
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)
$totalSheet.Cells.Item(2, 1).Value = 0

$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 0.22

# The row-insert shifted the old rows down but kept their old index values
# (0, 1) in column A; bump them to (1, 2) to match the renumbered list.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2

# Restore "总计" as the active/selected tab (matches original workbook state).
$totalSheet.Activate()
